$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 'sd'
$ws.Range("J2").Value = 'Statement-non-opinion'
$ws.Range("I5").Value = 'aa'
$ws.Range("J5").Value = 'Agree/Accept'
$ws.Range("I7").Value = 'sv'
$ws.Range("J7").Value = 'Statement-opinion'
$ws.Range("I10").Value = 'sd'
$ws.Range("J10").Value = 'Statement-non-opinion'
$ws.Range("I14").Value = 'sd'
$ws.Range("J14").Value = 'Statement-non-opinion'
$ws.Range("I21").Value = 'sd'
$ws.Range("J21").Value = 'Statement-non-opinion'
$ws.Range("I22").Value = 'b'
$ws.Range("J22").Value = 'Acknowledge (Backchannel)'
$ws.Range("I25").Value = 'sd'
$ws.Range("J25").Value = 'Statement-non-opinion'
$ws.Range("I31").Value = 'ba'
$ws.Range("J31").Value = 'Appreciation'
$ws.Range("I38").Value = 'ba'
$ws.Range("J38").Value = 'Appreciation'
$ws.Range("I40").Value = '%'
$ws.Range("J40").Value = 'Uninterpretable'
$ws.Range("I44").Value = 'aa'
$ws.Range("J44").Value = 'Agree/Accept'
$ws.Range("I46").Value = 'sv'
$ws.Range("J46").Value = 'Statement-opinion'
$ws.Range("I50").Value = 'sd'
$ws.Range("J50").Value = 'Statement-non-opinion'
$ws.Range("I72").Value = 'ba'
$ws.Range("J72").Value = 'Appreciation'
$ws.Range("I74").Value = 'sv'
$ws.Range("J74").Value = 'Statement-opinion'
$ws.Range("I75").Value = 'sd'
$ws.Range("J75").Value = 'Statement-non-opinion'
$ws.Range("I86").Value = 'sv'
$ws.Range("J86").Value = 'Statement-opinion'
$ws.Range("I98").Value = 'b'
$ws.Range("J98").Value = 'Acknowledge (Backchannel)'
$ws.Range("I99").Value = 'sv'
$ws.Range("J99").Value = 'Statement-opinion'
$ws.Range("I101").Value = 'sd'
$ws.Range("J101").Value = 'Statement-non-opinion'
$ws.Range("I103").Value = 'aa'
$ws.Range("J103").Value = 'Agree/Accept'
$ws.Range("I106").Value = 'sd'
$ws.Range("J106").Value = 'Statement-non-opinion'
$ws.Range("I134").Value = 'qy'
$ws.Range("J134").Value = 'Yes-No-Question'
$ws.Range("I153").Value = 'sv'
$ws.Range("J153").Value = 'Statement-opinion'
$ws.Range("I154").Value = 'sd'
$ws.Range("J154").Value = 'Statement-non-opinion'
$ws.Range("I162").Value = 'sd'
$ws.Range("J162").Value = 'Statement-non-opinion'
$ws.Range("I164").Value = 'sd'
$ws.Range("J164").Value = 'Statement-non-opinion'
$ws.Range("I166").Value = 'sd'
$ws.Range("J166").Value = 'Statement-non-opinion'
$ws.Range("I167").Value = 'aa'
$ws.Range("J167").Value = 'Agree/Accept'
$ws.Range("I168").Value = 'aa'
$ws.Range("J168").Value = 'Agree/Accept'
$ws.Range("I171").Value = 'aa'
$ws.Range("J171").Value = 'Agree/Accept'
$ws.Range("I174").Value = 'ba'
$ws.Range("J174").Value = 'Appreciation'
$ws.Range("I186").Value = 'sv'
$ws.Range("J186").Value = 'Statement-opinion'
$ws.Range("I188").Value = 'sv'
$ws.Range("J188").Value = 'Statement-opinion'
$ws.Range("I192").Value = 'sd'
$ws.Range("J192").Value = 'Statement-non-opinion'
$ws.Range("I203").Value = 'sv'
$ws.Range("J203").Value = 'Statement-opinion'
$ws.Range("I220").Value = 'sv'
$ws.Range("J220").Value = 'Statement-opinion'
$ws.Range("I227").Value = 'sd'
$ws.Range("J227").Value = 'Statement-non-opinion'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'sd'
$ws.Range("J236").Value = 'Statement-non-opinion'
$ws.Range("I238").Value = 'aa'
$ws.Range("J238").Value = 'Agree/Accept'
$ws.Range("I247").Value = '%'
$ws.Range("J247").Value = 'Uninterpretable'
$ws.Range("I251").Value = 'sd'
$ws.Range("J251").Value = 'Statement-non-opinion'
$ws.Range("I269").Value = 'b'
$ws.Range("J269").Value = 'Acknowledge (Backchannel)'
$ws.Range("I275").Value = '%'
$ws.Range("J275").Value = 'Uninterpretable'
$ws.Range("I303").Value = 'sd'
$ws.Range("J303").Value = 'Statement-non-opinion'
$ws.Range("I307").Value = 'aa'
$ws.Range("J307").Value = 'Agree/Accept'
$ws.Range("I313").Value = 'aa'
$ws.Range("J313").Value = 'Agree/Accept'
